# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 285
    4  = 10925
    5  = 10013
    8  = 711
    11 = 27
    13 = 9546
    15 = 2428
    17 = 76
    18 = 379
    19 = 10837
    20 = 10758
    21 = 7
    22 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
